$wb = $excel.ActiveWorkbook

# ---- "Metadata" sheet updates ----
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# New generation date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a real value
$meta.Range("B9").Value = "Alvearie Team"

# First "Contact" row becomes the new "Jurisdiction" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Second (duplicate) "Contact" row becomes the "Description" row
$meta.Range("A11").Value = "Description"
$meta.Range("B11").Value = "Age in months at the time of the event"

# Old "Description" row (now redundant) is removed; everything below shifts up
$meta.Rows.Item(12).Delete()

# ---- "Elements" sheet updates ----
$elements = $wb.Worksheets.Item("Elements")

# Root element Short / Definition updated to describe this extension
$elements.Range("K2").Value = "Snapshot Age in months"
$elements.Range("L2").Value = "Age in months at the time of the event"
